$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid = 1, Absent = 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count = 1, Real = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: Total Attendance Count = 1, Real = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: Absent = 1
$ws.Range("H6").Value = 1

# Row 7: Absent = 1
$ws.Range("H7").Value = 1

# Row 8: Absent = 1
$ws.Range("H8").Value = 1

# Row 9: Absent = 1
$ws.Range("H9").Value = 1

# Row 10: Absent = 1
$ws.Range("H10").Value = 1

# Row 11: Total Attendance Count = 1, Real = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1

# Row 12: Total Attendance Count = 1, Real = 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

# Row 13: Absent = 1
$ws.Range("H13").Value = 1

# Row 14: Absent = 1
$ws.Range("H14").Value = 1

# Row 15: Absent = 1
$ws.Range("H15").Value = 1

# Row 16: Absent = 1
$ws.Range("H16").Value = 1

# Row 17: Absent = 1
$ws.Range("H17").Value = 1

# Row 18: Absent = 1
$ws.Range("H18").Value = 1
